# Applies metadata updates to the "Metadata" worksheet of the CodeSystem workbook:
#  - Status changes from "draft" to "active"
#  - Experimental gets a value of "false"
#  - Date is refreshed to a newer timestamp
#  - Case Sensitive gets a value of "true"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B6").Value = "active"
$ws.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# "false"/"true" look like booleans to the engine's auto-type-detection, so a
# plain .Value assignment would store them as native booleans (t="b") rather
# than text. Routing them through a formula + values-only paste keeps them as
# literal text (t="s") without leaving stray formulas or quote-prefixed
# number formats behind.
$cExperimental = $ws.Range("B7")
$cExperimental.Formula = '="false"'
$cExperimental.Copy()
$cExperimental.PasteSpecial(-4163)

$cCaseSensitive = $ws.Range("B15")
$cCaseSensitive.Formula = '="true"'
$cCaseSensitive.Copy()
$cCaseSensitive.PasteSpecial(-4163)

$excel.CutCopyMode = 0
